$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-11 Saturday", "2025-10-12 Sunday"),
    @("371×3=1113", "169×6=1014"),
    @("107×4=428", "228×7=1596"),
    @("751×8=6008", "162×2=324"),
    @("514×5=2570", "264×3=792"),
    @("561×8=4488", "319×7=2233"),
    @("500×2=1000", "688×4=2752"),
    @("995×2=1990", "786×9=7074"),
    @("778×6=4668", "794×9=7146"),
    @("810×7=5670", "486×2=972"),
    @("963×6=5778", "875×2=1750"),
    @("261×9=2349", "743×2=1486"),
    @("958×4=3832", "214×6=1284"),
    @("127×9=1143", "811×8=6488"),
    @("773×3=2319", "824×5=4120"),
    @("747×3=2241", "158×7=1106"),
    @("351×9=3159", "164×2=328"),
    @("139×8=1112", "276×8=2208"),
    @("520×9=4680", "779×3=2337"),
    @("894×6=5364", "728×5=3640"),
    @("873×7=6111", "358×2=716"),
    @("566×3=1698", "397×5=1985"),
    @("536×2=1072", "162×9=1458"),
    @("465×7=3255", "106×7=742"),
    @("773×8=6184", "367×6=2202"),
    @("376×7=2632", "711×2=1422")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
